$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("A20").Value = "3/1, 1hr"
$ws.Range("B20").Value = "working on getting shiny app to work until today"

# Row 21
$ws.Range("A21").Value = "3/2, 30 mins"
$ws.Range("B21").Value = "testing covid overtime graph"

# Copy style from row 19 (A19/B19) to the new rows 20/21 so formatting matches
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 20 should have the same wrapped row height as row 19, row 21 stays default
$ws.Rows.Item(20).RowHeight = $ws.Rows.Item(19).RowHeight

# Update selection to match the new active cell position
$ws.Range("B22").Select()
